$d = $word.ActiveDocument

function Remove-ParagraphByText($searchText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $true, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $para = $rng.Paragraphs(1)
        $para.Range.Delete()
    }
}

# Remove the "Develop wireframes and mockups ..." paragraph entirely
# (it sat directly under "User Interface Design:").
Remove-ParagraphByText("Develop wireframes and*experience.")

# Remove the "Define the data models ..." paragraph entirely
# (it sat directly under "Data Models:").
Remove-ParagraphByText("Define the data models*user information.")

# Remove the "Detail the security measures ..." paragraph entirely
# (it sat directly under "Security Plan:").
Remove-ParagraphByText("Detail the security measures*at rest.")
